# Auto-generated Excel COM-interop script
# Applies numeric value updates to the Leve profit tables across all 8 sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 76
$ws.Range("H76").Value = 3068.3076
$ws.Range("I76").Value = 3002.7273
$ws.Range("J76").Value = 3153.1765
$ws.Range("K76").Value = 3002.7273
$ws.Range("L76").Value = 3153.1765
$ws.Range("M76").Value = -2687.7273
$ws.Range("N76").Value = -3783.1765
# Row 79
$ws.Range("H79").Value = 3068.3076
$ws.Range("I79").Value = 3002.7273
$ws.Range("J79").Value = 3153.1765
$ws.Range("K79").Value = 3002.7273
$ws.Range("L79").Value = 3153.1765
$ws.Range("M79").Value = -1910.7273
$ws.Range("N79").Value = -5337.1765
# Row 115
$ws.Range("H115").Value = 1942.8462
$ws.Range("I115").Value = 1942.8462
$ws.Range("K115").Value = 5828.5386
$ws.Range("M115").Value = -4261.5386

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7962.1797
$ws.Range("I32").Value = 6824.811
$ws.Range("K32").Value = 6824.811
$ws.Range("M32").Value = -6537.811
# Row 132
$ws.Range("H132").Value = 1632.8214
$ws.Range("I132").Value = 942.86365
$ws.Range("K132").Value = 2828.59095
$ws.Range("M132").Value = -298.5909499999998

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 1857.8182
$ws.Range("I107").Value = 1500.5714
$ws.Range("J107").Value = 2483
$ws.Range("K107").Value = 1500.5714
$ws.Range("L107").Value = 2483
$ws.Range("M107").Value = 419.4286
$ws.Range("N107").Value = -6323
# Row 134
$ws.Range("H134").Value = 2341.1392
$ws.Range("I134").Value = 1508.9783
$ws.Range("K134").Value = 4526.9349
$ws.Range("M134").Value = -1991.9349

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1191.6086
$ws.Range("I16").Value = 1074.9474
$ws.Range("J16").Value = 1745.75
$ws.Range("K16").Value = 1074.9474
$ws.Range("L16").Value = 1745.75
$ws.Range("M16").Value = -787.9474
$ws.Range("N16").Value = -2319.75
# Row 31
$ws.Range("H31").Value = 3053.96
$ws.Range("I31").Value = 1726.5807
$ws.Range("J31").Value = 3650.3188
$ws.Range("K31").Value = 1726.5807
$ws.Range("L31").Value = 3650.3188
$ws.Range("M31").Value = -1431.5807
$ws.Range("N31").Value = -4240.3188
# Row 34
$ws.Range("H34").Value = 3053.96
$ws.Range("I34").Value = 1726.5807
$ws.Range("J34").Value = 3650.3188
$ws.Range("K34").Value = 1726.5807
$ws.Range("L34").Value = 3650.3188
$ws.Range("M34").Value = -1524.5807
$ws.Range("N34").Value = -4054.3188
# Row 86
$ws.Range("H86").Value = 4371.385
$ws.Range("I86").Value = 4166.273
$ws.Range("J86").Value = 5499.5
$ws.Range("K86").Value = 4166.273
$ws.Range("L86").Value = 5499.5
$ws.Range("M86").Value = -3043.273
$ws.Range("N86").Value = -7745.5
# Row 89
$ws.Range("H89").Value = 4371.385
$ws.Range("I89").Value = 4166.273
$ws.Range("J89").Value = 5499.5
$ws.Range("K89").Value = 20831.365
$ws.Range("L89").Value = 27497.5
$ws.Range("M89").Value = -15215.365
$ws.Range("N89").Value = -38729.5
# Row 113
$ws.Range("H113").Value = 1191.6086
$ws.Range("I113").Value = 1074.9474
$ws.Range("J113").Value = 1745.75
$ws.Range("K113").Value = 1074.9474
$ws.Range("L113").Value = 1745.75
$ws.Range("M113").Value = 1095.0526
$ws.Range("N113").Value = -6085.75
# Row 134
$ws.Range("H134").Value = 71377.60000000001
$ws.Range("I134").Value = 1073.4615
$ws.Range("K134").Value = 3220.3845
$ws.Range("M134").Value = -685.3844999999997

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 17
$ws.Range("H17").Value = 1025
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 1200
$ws.Range("K17").Value = 1500
$ws.Range("L17").Value = 3600
$ws.Range("M17").Value = -1331
$ws.Range("N17").Value = -3938
# Row 23
$ws.Range("H23").Value = 593.6667
$ws.Range("I23").Value = 1399
$ws.Range("J23").Value = 493
$ws.Range("K23").Value = 4197
$ws.Range("L23").Value = 1479
$ws.Range("M23").Value = -3962
$ws.Range("N23").Value = -1949
# Row 33
$ws.Range("H33").Value = 15703070
$ws.Range("I33").Value = 106.333336
$ws.Range("J33").Value = 24268324
$ws.Range("K33").Value = 638.000016
$ws.Range("L33").Value = 145609944
$ws.Range("M33").Value = -355.000016
$ws.Range("N33").Value = -145610510
# Row 98
$ws.Range("H98").Value = 1107.4
$ws.Range("I98").Value = 400
$ws.Range("J98").Value = 1186
$ws.Range("K98").Value = 1200
$ws.Range("L98").Value = 3558
$ws.Range("M98").Value = 298
$ws.Range("N98").Value = -6554
# Row 101
$ws.Range("H101").Value = 6000
$ws.Range("J101").Value = 6000
$ws.Range("L101").Value = 18000
$ws.Range("N101").Value = -22868
# Row 102
$ws.Range("H102").Value = 6988.4443
$ws.Range("J102").Value = 6612.125
$ws.Range("L102").Value = 19836.375
$ws.Range("N102").Value = -24704.375
# Row 104
$ws.Range("H104").Value = 1570
$ws.Range("J104").Value = 1587.5
$ws.Range("L104").Value = 4762.5
$ws.Range("N104").Value = -10004.5
# Row 105
$ws.Range("H105").Value = 129058.5
$ws.Range("J105").Value = 147138.28
$ws.Range("L105").Value = 441414.84
$ws.Range("N105").Value = -446656.84
# Row 131
$ws.Range("H131").Value = 854.51
$ws.Range("I131").Value = 486.5
$ws.Range("K131").Value = 1459.5
$ws.Range("M131").Value = 3580.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5195.8335
$ws.Range("I70").Value = 5236.8423
$ws.Range("J70").Value = 5040
$ws.Range("K70").Value = 5236.8423
$ws.Range("L70").Value = 5040
$ws.Range("M70").Value = -4966.8423
$ws.Range("N70").Value = -5580
# Row 73
$ws.Range("H73").Value = 5195.8335
$ws.Range("I73").Value = 5236.8423
$ws.Range("J73").Value = 5040
$ws.Range("K73").Value = 5236.8423
$ws.Range("L73").Value = 5040
$ws.Range("M73").Value = -4300.8423
$ws.Range("N73").Value = -6912
# Row 97
$ws.Range("H97").Value = 4105.593
$ws.Range("I97").Value = 2326.6667
$ws.Range("J97").Value = 18337
$ws.Range("K97").Value = 2326.6667
$ws.Range("L97").Value = 18337
$ws.Range("M97").Value = -1830.6667
$ws.Range("N97").Value = -19329

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 3044.8333
$ws.Range("I61").Value = 3191.8572
$ws.Range("J61").Value = 2839
$ws.Range("K61").Value = 3191.8572
$ws.Range("L61").Value = 2839
$ws.Range("M61").Value = -2989.8572
$ws.Range("N61").Value = -3243
# Row 100
$ws.Range("H100").Value = 2767.5
$ws.Range("I100").Value = 2767.5
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2767.5
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -2226.5
# Row 113
$ws.Range("H113").Value = 3044.8333
$ws.Range("I113").Value = 3191.8572
$ws.Range("J113").Value = 2839
$ws.Range("K113").Value = 3191.8572
$ws.Range("L113").Value = 2839
$ws.Range("M113").Value = -1021.8572
$ws.Range("N113").Value = -7179
# Row 136
$ws.Range("H136").Value = 1834.6364
$ws.Range("I136").Value = 1433.1305
$ws.Range("J136").Value = 2758.1
$ws.Range("K136").Value = 4299.3915
$ws.Range("L136").Value = 8274.299999999999
$ws.Range("M136").Value = -1749.3915
$ws.Range("N136").Value = -13374.3

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 677.4286
$ws.Range("I100").Value = 440.66666
$ws.Range("J100").Value = 855
$ws.Range("K100").Value = 881.33332
$ws.Range("L100").Value = 1710
$ws.Range("M100").Value = -340.33332
$ws.Range("N100").Value = -2792
# Row 101
$ws.Range("H101").Value = 20950
$ws.Range("J101").Value = 20950
$ws.Range("L101").Value = 20950
$ws.Range("N101").Value = -27440
